# Auto-generated PowerShell COM-interop script
# Applies numeric value updates to the Ultros_Profits.xlsx market-price sheets
# as captured by the authoritative OOXML diff (commit: "chore: update Sheets via scheduled runner").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 397.77777
$ws.Range("I9").Value = 270.25
$ws.Range("J9").Value = 499.8
$ws.Range("K9").Value = 270.25
$ws.Range("L9").Value = 499.8
$ws.Range("M9").Value = -101.25
$ws.Range("N9").Value = -837.8

$ws.Range("H15").Value = 2979.4119
$ws.Range("I15").Value = 2979.4119
$ws.Range("K15").Value = 8938.235700000001
$ws.Range("M15").Value = -8769.235700000001

$ws.Range("H32").Value = 9094041
$ws.Range("J32").Value = 10003355
$ws.Range("L32").Value = 10003355
$ws.Range("N32").Value = -10004007

$ws.Range("H43").Value = 5750
$ws.Range("I43").Value = 5000
$ws.Range("J43").Value = 6500
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 6500
$ws.Range("M43").Value = -4931
$ws.Range("N43").Value = -6638

$ws.Range("H92").Value = 2485.125
$ws.Range("I92").Value = 2580.2307
$ws.Range("K92").Value = 2580.2307
$ws.Range("M92").Value = -1332.2307

$ws.Range("H96").Value = 1373.3334
$ws.Range("I96").Value = 1060
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 3180
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -1807
$ws.Range("N96").Value = -8746

$ws.Range("H125").Value = 14709.2
$ws.Range("I125").Value = 21253.8
$ws.Range("K125").Value = 191284.2
$ws.Range("M125").Value = -188824.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 14269.741
$ws.Range("I2").Value = 17794.3
$ws.Range("K2").Value = 17794.3
$ws.Range("M2").Value = -17681.3

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H74").Value = 1580.75
$ws.Range("I74").Value = 1580.75
$ws.Range("K74").Value = 1580.75
$ws.Range("M74").Value = -706.75

$ws.Range("H77").Value = 1580.75
$ws.Range("I77").Value = 1580.75
$ws.Range("K77").Value = 7903.75
$ws.Range("M77").Value = -3535.75

$ws.Range("H116").Value = 14269.741
$ws.Range("I116").Value = 17794.3
$ws.Range("K116").Value = 17794.3
$ws.Range("M116").Value = -15500.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 14269.741
$ws.Range("I3").Value = 17794.3
$ws.Range("K3").Value = 17794.3
$ws.Range("M3").Value = -17680.3

$ws.Range("H134").Value = 2555.2222
$ws.Range("I134").Value = 2063.4546
$ws.Range("J134").Value = 3328
$ws.Range("K134").Value = 6190.3638
$ws.Range("L134").Value = 9984
$ws.Range("M134").Value = -3655.3638
$ws.Range("N134").Value = -15054

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3599
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 3748.75
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 3748.75
$ws.Range("M58").Value = -2797
$ws.Range("N58").Value = -4154.75

$ws.Range("H94").Value = 2099.5715
$ws.Range("J94").Value = 2216
$ws.Range("L94").Value = 2216
$ws.Range("N94").Value = -3118

$ws.Range("H132").Value = 2881.6667
$ws.Range("I132").Value = 2119.7334
$ws.Range("J132").Value = 6691.3335
$ws.Range("K132").Value = 6359.2002
$ws.Range("L132").Value = 20074.0005
$ws.Range("M132").Value = -3829.2002
$ws.Range("N132").Value = -25134.0005

$ws.Range("H136").Value = 3599
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 3748.75
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 11246.25
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -16346.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 684.3
$ws.Range("I8").Value = 684.3
$ws.Range("K8").Value = 2052.9
$ws.Range("M8").Value = -1913.9

$ws.Range("H87").Value = 1200
$ws.Range("I87").Value = 1200
$ws.Range("K87").Value = 3600
$ws.Range("M87").Value = -2352

$ws.Range("H90").Value = 1200
$ws.Range("I90").Value = 1200
$ws.Range("K90").Value = 10800
$ws.Range("M90").Value = -4560

$ws.Range("H98").Value = 2179.2856
$ws.Range("I98").Value = 2373.2856
$ws.Range("J98").Value = 1985.2858
$ws.Range("K98").Value = 7119.8568
$ws.Range("L98").Value = 5955.857400000001
$ws.Range("M98").Value = -5621.8568
$ws.Range("N98").Value = -8951.857400000001

$ws.Range("H121").Value = 1021.9
$ws.Range("I121").Value = 537.4
$ws.Range("J121").Value = 1506.4
$ws.Range("K121").Value = 1612.2
$ws.Range("L121").Value = 4519.200000000001
$ws.Range("M121").Value = -302.1999999999998
$ws.Range("N121").Value = -7139.200000000001

$ws.Range("H131").Value = 2268.3333
$ws.Range("I131").Value = 1916.5555
$ws.Range("J131").Value = 2796
$ws.Range("K131").Value = 5749.666499999999
$ws.Range("L131").Value = 8388
$ws.Range("M131").Value = -709.6664999999994
$ws.Range("N131").Value = -18468

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 33799.2
$ws.Range("J52").Value = 34749
$ws.Range("L52").Value = 34749
$ws.Range("N52").Value = -35267

$ws.Range("H70").Value = 69225.11
$ws.Range("I70").Value = 85361.78999999999
$ws.Range("K70").Value = 85361.78999999999
$ws.Range("M70").Value = -85091.78999999999

$ws.Range("H73").Value = 69225.11
$ws.Range("I73").Value = 85361.78999999999
$ws.Range("K73").Value = 85361.78999999999
$ws.Range("M73").Value = -84425.78999999999

$ws.Range("H102").Value = 3668.92
$ws.Range("I102").Value = 2895.9443
$ws.Range("K102").Value = 2895.9443
$ws.Range("M102").Value = -1273.9443

$ws.Range("H107").Value = 323.86957
$ws.Range("I107").Value = 356.06668
$ws.Range("J107").Value = 263.5
$ws.Range("K107").Value = 356.06668
$ws.Range("L107").Value = 263.5
$ws.Range("M107").Value = 1563.93332
$ws.Range("N107").Value = -4103.5

$ws.Range("H113").Value = 9678.786
$ws.Range("J113").Value = 12998.6
$ws.Range("L113").Value = 12998.6
$ws.Range("N113").Value = -17338.6

$ws.Range("H123").Value = 35599
$ws.Range("J123").Value = 35599
$ws.Range("L123").Value = 35599
$ws.Range("N123").Value = -40499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4072.0908
$ws.Range("I7").Value = 2458.6
$ws.Range("J7").Value = 5416.6665
$ws.Range("K7").Value = 2458.6
$ws.Range("L7").Value = 5416.6665
$ws.Range("M7").Value = -2346.6
$ws.Range("N7").Value = -5640.6665

$ws.Range("H40").Value = 2651.75
$ws.Range("I40").Value = 2651.75
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2651.75
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2515.75
$ws.Range("N40").ClearContents()

$ws.Range("H126").Value = 4072.0908
$ws.Range("I126").Value = 2458.6
$ws.Range("J126").Value = 5416.6665
$ws.Range("K126").Value = 7375.799999999999
$ws.Range("L126").Value = 16249.9995
$ws.Range("M126").Value = -4905.799999999999
$ws.Range("N126").Value = -21189.9995

$ws.Range("H136").Value = 4364.1113
$ws.Range("I136").Value = 3886.7273
$ws.Range("K136").Value = 11660.1819
$ws.Range("M136").Value = -9110.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 501.5
$ws.Range("I107").Value = 501.5
$ws.Range("K107").Value = 1504.5
$ws.Range("M107").Value = 415.5

$ws.Range("H122").Value = 2158.4
$ws.Range("J122").Value = 2266
$ws.Range("L122").Value = 6798
$ws.Range("N122").Value = -11698

$ws.Range("H132").Value = 4704.5884
$ws.Range("I132").Value = 4405.6665
$ws.Range("K132").Value = 13216.9995
$ws.Range("M132").Value = -10686.9995

